$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows continuing the existing time series (update through 21 marzo / new rows 230-233)
$rows = @(
    @{ Row = 230; A = 44304; B = 4; C = 27; D = 174.3848091455144 },
    @{ Row = 231; A = 44305; B = 5; C = 27; D = 174.3848091455144 },
    @{ Row = 232; A = 44306; B = 3; C = 26; D = 167.9261125104954 },
    @{ Row = 233; A = 44307; B = 2; C = 27; D = 174.3848091455144 }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Replicate the date-column formatting (border/alignment/number-format) from
    # the last existing data row so the new cells carry the same style index.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
}

$excel.CutCopyMode = 0
